$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "6eca91e3ee06c8bff1654dba0f3e21432c104d43876eefbf6b48c7d0e4cea770"
$ws.Range("D3").Value = "262911b1308e8175e3ad988887fea8be3ba8bcae08e428719721267456123163"
$ws.Range("D4").Value = "c59557639441bcd0ed68a07b4555f7f7bb7702e6b97677d67a982ab4400884fe"
$ws.Range("D5").Value = "321c69c926111abfa23c2d89884607c968af7aab8e3d7eabbf42bd7cad26b198"
$ws.Range("D6").Value = "bdf901631cafc7894c25c38edd8e8c1d1afab21fe2419292653a17972ee558f4"
